$wb = $excel.ActiveWorkbook

# 1. Rename first sheet "Tests" -> "General tests"
$wsGeneral = $wb.Worksheets.Item(1)
$wsGeneral.Name = "General tests"

# 2. Add new sheet right after "General tests", named "Twitter tests"
$ws = $wb.Worksheets.Add($null, $wsGeneral)
$ws.Name = "Twitter tests"

$ws.Range('A1').Value = 'Tweeter'
$ws.Range('B1').Value = 'Nr words'
$ws.Range('F1').Value = 'SKSS PERS'
$ws.Range('E1').Value = 'CKS PERS'
$ws.Range('C1').Value = 'CKS NL'
$ws.Range('D1').Value = 'SKSS NL'
$ws.Range('A2').Value = 'artbysophia'
$ws.Range('A3').Value = 'baspaternotte'
$ws.Range('A4').Value = 'marktwain2'
$ws.Range('A5').Value = 'pinaatje'
$ws.Range('A6').Value = 'pienbetuwe'
$ws.Range('A7').Value = 'rider_ot_storm'
$ws.Range('A8').Value = 'mrsmartine'
$ws.Range('A9').Value = 'chrisklomp'
$ws.Range('A10').Value = 'leolewin'
$ws.Range('A11').Value = 'contentgirl'
$ws.Range('A12').Value = 'amadeusivan'
$ws.Range('A13').Value = 'ongerijmd'
$ws.Range('A14').Value = 'umarebru'
$ws.Range('A15').Value = 'a_mieke'
$ws.Range('A16').Value = 'steephsel'
$ws.Range('A17').Value = 'brechtjedeleij'
$ws.Range('I1').Value = 'Opmerkingen'
$ws.Range('I12').Value = 'mix met en'
$ws.Range('A18').Value = 'eetschrijver'
$ws.Range('A20').Value = 'klapster'
$ws.Range('A21').Value = 'goedemorgenman'
$ws.Range('A22').Value = 'walterhoekstra'
$ws.Range('A23').Value = 'jasmijn02'
$ws.Range('A24').Value = 'miekeinc'
$ws.Range('A25').Value = 'sredlums'
$ws.Range('A26').Value = 'aldith_hunkar'
$ws.Range('A27').Value = 'tien020'
$ws.Range('A28').Value = 'karinwinters'
$ws.Range('A29').Value = 'johnschop'
$ws.Range('A30').Value = 'lobdozer'
$ws.Range('A31').Value = 'theollieworks'
$ws.Range('A32').Value = 'wup5'
$ws.Range('A33').Value = 'jennekepenneke'
$ws.Range('A34').Value = 'rebelsnotes'
$ws.Range('A35').Value = 'puberdochters'
$ws.Range('A36').Value = 'knotsbots'
$ws.Range('A37').Value = 'dennismons'
$ws.Range('A38').Value = 'fluist3r'
$ws.Range('A39').Value = 'mariannecramer'
$ws.Range('A40').Value = 'rjvanhouten'
$ws.Range('A41').Value = 'superjan'
$ws.Range('A42').Value = 'titchener'
$ws.Range('A43').Value = 'anniebbarks'
$ws.Range('A44').Value = 'fred3012'
$ws.Range('A45').Value = 'politicus1'
$ws.Range('A46').Value = 'peterstafleu'
$ws.Range('A47').Value = 'jettyvanrooy'
$ws.Range('A48').Value = 'mariannezw'
$ws.Range('A49').Value = 'jochemgeerdink'
$ws.Range('A50').Value = 'nabilfeki'
$ws.Range('A51').Value = 'kos_'
$ws.Range('A19').Value = 'esther_305'
$ws.Range('G1').Value = 'CKS COMMUNITY'
$ws.Range('H1').Value = 'SKSS COMMUNITY'
$ws.Range('M1').Value = 'SKKS OTHER'
$ws.Range('L1').Value = 'CKS OTHER'
$ws.Range('B2').Value = 477584
$ws.Range('C2').Value = 21
$ws.Range('D2').Value = 21
$ws.Range('E2').Value = 27
$ws.Range('F2').Value = 31
$ws.Range('G2').Value = 36
$ws.Range('H2').Value = 36
$ws.Range('L2').Value = 33
$ws.Range('M2').Value = 37
$ws.Range('N2').Value = 'pienbetuwe'
$ws.Range('B3').Value = 324268
$ws.Range('C3').Value = 16
$ws.Range('D3').Value = 16
$ws.Range('E3').Value = 26
$ws.Range('F3').Value = 28
$ws.Range('G3').Value = 31
$ws.Range('H3').Value = 34
$ws.Range('L3').Value = 28
$ws.Range('M3').Value = 31
$ws.Range('N3').Value = 'artbysophia'
$ws.Range('B4').Value = 334816
$ws.Range('C4').Value = 20
$ws.Range('D4').Value = 20
$ws.Range('E4').Value = 33
$ws.Range('F4').Value = 37
$ws.Range('G4').Value = 35
$ws.Range('H4').Value = 39
$ws.Range('L4').Value = 33
$ws.Range('M4').Value = 37
$ws.Range('N4').Value = 'pinaatje'
$ws.Range('B5').Value = 297352
$ws.Range('C5').Value = 17
$ws.Range('D5').Value = 17
$ws.Range('E5').Value = 30
$ws.Range('F5').Value = 33
$ws.Range('G5').Value = 34
$ws.Range('H5').Value = 36
$ws.Range('L5').Value = 33
$ws.Range('M5').Value = 37
$ws.Range('N5').Value = 'baspaternotte'
$ws.Range('B6').Value = 253684
$ws.Range('C6').Value = 16
$ws.Range('D6').Value = 17
$ws.Range('E6').Value = 26
$ws.Range('F6').Value = 30
$ws.Range('G6').Value = 31
$ws.Range('H6').Value = 36
$ws.Range('L6').Value = 30
$ws.Range('M6').Value = 36
$ws.Range('N6').Value = 'marktwain2'
$ws.Range('B7').Value = 291008
$ws.Range('C7').Value = 14
$ws.Range('D7').Value = 14
$ws.Range('E7').Value = 38
$ws.Range('F7').Value = 44
$ws.Range('B8').Value = 330960
$ws.Range('C8').Value = 19
$ws.Range('D8').Value = 19
$ws.Range('E8').Value = 25
$ws.Range('F8').Value = 28
$ws.Range('B9').Value = 316428
$ws.Range('C9').Value = 20
$ws.Range('D9').Value = 20
$ws.Range('E9').Value = 24
$ws.Range('F9').Value = 28
$ws.Range('B10').Value = 286320
$ws.Range('C10').Value = 22
$ws.Range('D10').Value = 22
$ws.Range('E10').Value = 34
$ws.Range('F10').Value = 38
$ws.Range('B11').Value = 272952
$ws.Range('C11').Value = 17
$ws.Range('D11').Value = 17
$ws.Range('E11').Value = 22
$ws.Range('F11').Value = 26
$ws.Range('B12').Value = 258364
$ws.Range('C12').Value = 9
$ws.Range('D12').Value = 9
$ws.Range('E12').Value = 40
$ws.Range('F12').Value = 48
$ws.Range('B13').Value = 218380
$ws.Range('C13').Value = 15
$ws.Range('D13').Value = 26
$ws.Range('E13').Value = 29
$ws.Range('F13').Value = 35
$ws.Range('B14').Value = 297144
$ws.Range('C14').Value = 19
$ws.Range('D14').Value = 19
$ws.Range('E14').Value = 24
$ws.Range('F14').Value = 27
$ws.Range('B15').Value = 260072
$ws.Range('C15').Value = 18
$ws.Range('D15').Value = 18
$ws.Range('E15').Value = 27
$ws.Range('F15').Value = 31
$ws.Range('B16').Value = 252460
$ws.Range('C16').Value = 17
$ws.Range('D16').Value = 17
$ws.Range('E16').Value = 24
$ws.Range('F16').Value = 27
$ws.Range('B17').Value = 201542
$ws.Range('C17').Value = 17
$ws.Range('D17').Value = 17
$ws.Range('E17').Value = 27
$ws.Range('F17').Value = 30
$ws.Range('B18').Value = 276976
$ws.Range('C18').Value = 20
$ws.Range('D18').Value = 20
$ws.Range('E18').Value = 30
$ws.Range('F18').Value = 33
$ws.Range('B19').Value = 256432
$ws.Range('C19').Value = 18
$ws.Range('D19').Value = 18
$ws.Range('E19').Value = 23
$ws.Range('F19').Value = 26
$ws.Range('B20').Value = 249304
$ws.Range('C20').Value = 19
$ws.Range('D20').Value = 19
$ws.Range('E20').Value = 25
$ws.Range('F20').Value = 28
$ws.Range('E21').Value = 50
$ws.Range('F21').Value = 53
$ws.Range('E23').Value = 30
$ws.Range('F23').Value = 34
$ws.Range('E24').Value = 21
$ws.Range('F24').Value = 24
$ws.Range('E25').Value = 25
$ws.Range('F25').Value = 26
$ws.Range('E26').Value = 25
$ws.Range('F26').Value = 29
$ws.Range('E27').Value = 26
$ws.Range('F27').Value = 30

# Column widths for Twitter tests sheet (closest achievable under engine's 1/6-char rounding)
$ws.Columns.Item(1).ColumnWidth = 15
$ws.Columns.Item(2).ColumnWidth = 13.666666666666666

# Selection on Twitter tests sheet
$ws.Range("F20").Select()

# 3. Training dirs sheet: add "Nr words" column
$wsTrain = $wb.Worksheets.Item("Training dirs")
$wsTrain.Range("E1").Value = "Nr words"
$wsTrain.Range("E4").Value = 55212868
$wsTrain.Activate()
$wsTrain.Range("A7").Select()

# 4. Testfiles sheet: update selection only
$wsTestfiles = $wb.Worksheets.Item("Testfiles")
$wsTestfiles.Activate()
$wsTestfiles.Range("B4").Select()

# 5. Re-activate Twitter tests sheet so it is the selected/visible tab
$ws.Activate()
